# Atualização de bases das ligas, do dia: 12-06-2024 às 23:38
#
# Several match rows got re-sorted/re-ordered in the source data. Column A
# (the running row index) stays put, but all other columns (B..AD) for the
# affected row pairs are swapped with one another.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param(
        [int]$Row1,
        [int]$Row2,
        [int]$FirstCol,
        [int]$LastCol
    )

    for ($col = $FirstCol; $col -le $LastCol; $col++) {
        $cell1 = $ws.Cells.Item($Row1, $col)
        $cell2 = $ws.Cells.Item($Row2, $col)

        $val1 = $cell1.Value2
        $val2 = $cell2.Value2

        $cell1.Value2 = $val2
        $cell2.Value2 = $val1
    }
}

# Swap pairs of rows whose B..AD data was interchanged.
# NOTE: this runtime's PowerShell parser does not reliably bind named
# (-Param value) arguments or parameter defaults, so all arguments are
# passed positionally and explicitly.
Swap-RowData 213 217 2 30
Swap-RowData 214 216 2 30
Swap-RowData 238 239 2 30
Swap-RowData 240 241 2 30
